# Apply the "Sin Deteccion" row-2 removal and refreshed "Resumen" metrics.

$wb = $excel.ActiveWorkbook

# --- Sheet "Sin Deteccion": drop the stray data row (A2=22, B2/C2 blank) ---
$wsSinDeteccion = $wb.Worksheets.Item("Sin Deteccion")
$wsSinDeteccion.Rows.Item(2).Delete()

# --- Sheet "Resumen": refresh the computed metrics ---
$wsResumen = $wb.Worksheets.Item("Resumen")
$wsResumen.Range("B2").Value = 0
$wsResumen.Range("B4").Value = 0

# B5 holds a literal text percentage (not a real number) - force text
# formatting before the assignment so Excel doesn't "smart" convert the
# string into a numeric percent value, then restore the default style so
# the cell format itself is left untouched.
$wsResumen.Range("B5").NumberFormat = "@"
$wsResumen.Range("B5").Value = "0%"
$wsResumen.Range("B5").Style = "Normal"

$wsResumen.Range("B6").Value = "2026-02-06 05:27:12"
